# Backlog.xlsx update — "subindo versão final site"
#
# Summary of the edit being reproduced:
#  - Most backlog items move to Status = "Finalizado" (project wrap-up).
#  - Three items get a Description (column B) filled in.
#  - The obsolete "Fazer Dashboard em página HTML" row is removed
#    (its row shifts everything below it up by one).
#  - A new row is appended at the bottom: "Fazer diagrama técnico".
#  - Column B (Descrição) gains center+border formatting across the table.
#  - The active selection ends on F18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 2) Status (column F) updates -> "Finalizado" for the rows that were
#    wrapped up. Each status has its own look (font color + fill) that
#    already exists elsewhere in the sheet, so copy the format from an
#    existing "Finalizado" cell (F5) and then set the value.
# ---------------------------------------------------------------------
function Set-Status($row, $status, $templateRow) {
    $ws.Range("F$templateRow").Copy() | Out-Null
    $ws.Range("F$row").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range("F$row").Value = $status
}

foreach ($row in @(3, 6, 7, 8, 9, 10, 11, 13, 14)) {
    Set-Status $row "Finalizado" 5
}

# ---------------------------------------------------------------------
# 3) Fill in the new descriptions (column B) for a few rows.
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "Criar as rotas, fazendo as conexões do front end com o banco de dados"

# ---------------------------------------------------------------------
# 4) Remove the obsolete row ("Fazer Dashboard em página HTML", row 15).
#    Everything below shifts up by one row.
# ---------------------------------------------------------------------
$ws.Range("A15").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------
# 5) After the shift, rows 15 (was 16) and 16 (was 17) finish up, plus
#    row 17 (was 18) status, and the rest of the statuses that still
#    needed to flip to "Finalizado".
# ---------------------------------------------------------------------
$ws.Range("B15").Value = "Fazer o sistema de interação do site com o usuário"
Set-Status 15 "Finalizado" 5

$ws.Range("B16").Value = "Configurar i banco de dados MySQL para rodar na VM virtual box"

Set-Status 17 "Finalizado" 5
Set-Status 20 "Finalizado" 5

# ---------------------------------------------------------------------
# 6) New row 21: "Fazer diagrama técnico".
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Fazer diagrama técnico"
$ws.Range("A21").Borders.Item(7).LineStyle = 1
$ws.Range("A21").Borders.Item(10).LineStyle = 1

$ws.Range("B21").Value = "Fazer o diagrama de solução técnico com as etapas do serviço realizado"
$ws.Range("B21").HorizontalAlignment = -4108

$ws.Range("C20").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("C21").Value = "Importante"

$ws.Range("F20").Copy() | Out-Null
$ws.Range("F21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F21").Value = "Finalizado"

# ---------------------------------------------------------------------
# 7) Column B (Descrição), rows 3-20: add the border+center formatting
#    that the rest of the data columns already have (copy format from
#    C3, which already carries that exact style). Done last, after the
#    row shift, so it lands on the final row numbering.
# ---------------------------------------------------------------------
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B3:B20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 8) Update the active selection to match.
# ---------------------------------------------------------------------
$ws.Range("F18").Select() | Out-Null
